# Rough outline of a schema: rebuild the workbook's sheets to reflect
# the pfhub schema classes (Benchmark, BenchmarkResult, Contributor,
# File, Implementation) -- each present twice (plain + "1" suffixed).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Target sheets, in final tab order, with their header rows.
# The first 5 reuse/rename the 6 existing sheets (positions 1-6);
# the remaining 4 are appended as brand-new sheets.
# ---------------------------------------------------------------------
$targets = @(
    @{ Name = "Benchmark";         Headers = @("id","name","version") },
    @{ Name = "BenchmarkResult";   Headers = @("id","name","contributors","benchmark","framework","hardware","implementation") },
    @{ Name = "Contributor";       Headers = @("id","name","handle","affiliation","email") },
    @{ Name = "File";              Headers = @("id","name","format") },
    @{ Name = "Implementation";    Headers = @("id","name","repository") },
    @{ Name = "Benchmark1";        Headers = @("id","name","version") },
    @{ Name = "BenchmarkResult1";  Headers = @("id","name","contributors","benchmark","framework","hardware","implementation") },
    @{ Name = "Contributor1";      Headers = @("id","name","handle","affiliation","email") },
    @{ Name = "File1";             Headers = @("id","name","format") },
    @{ Name = "Implementation1";   Headers = @("id","name","repository") }
)

$existingCount = $wb.Worksheets.Count

for ($i = 0; $i -lt $targets.Count; $i++) {
    $idx = $i + 1
    $target = $targets[$i]

    if ($idx -le $existingCount) {
        $ws = $wb.Worksheets.Item($idx)
    } else {
        $last = $wb.Worksheets.Item($wb.Worksheets.Count)
        $ws = $wb.Worksheets.Add($null, $last)
    }

    # Wipe any previous content/validation so only the new header survives.
    $ws.Cells.Clear()
    $ws.Cells.Validation.Delete()

    $ws.Name = $target.Name

    $headers = $target.Headers
    for ($c = 0; $c -lt $headers.Count; $c++) {
        $cell = $ws.Cells.Item(1, $c + 1)
        $cell.Value = $headers[$c]
    }
}

# The last sheet (Implementation1) ends up active, matching the source
# workbook having its newest sheet selected.
$wb.Worksheets.Item($targets.Count).Activate()
